$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('Símbolo','First','Follow'),
    @('program ','init','$'),
    @('program-suffix','id, if, do, read, write','$'),
    @('decl-assign',':=, ",", is','stop'),
    @('other-stmt','if, do, read, write','stop'),
    @('stmt-prime','if, do, read, write',';'),
    @('decl-list ','id, if, do, read, write','id, read, write, if, do'),
    @('ident-list ','",", λ','is'),
    @('type ','integer, string',';'),
    @('stmt-list ','id, read, write, if, do','end, while'),
    @('stmt-list-tail ','id, read, write, if, do, λ','stop, end, while'),
    @('stmt ','id, read, write, if, do',';'),
    @('assign-stmt ','id',';'),
    @('if-stmt ','if',';'),
    @('if-suffix ','else, λ',';'),
    @('condition ','id, num, string, (, not, -',')'),
    @('do-stmt ','do',';'),
    @('do-suffix ','while',';'),
    @('read-stmt ','read',';'),
    @('write-stmt ','write',';'),
    @('writable ','id, num, string, (, not, -',')'),
    @('expression ','id, num, string, (, not, -',')'),
    @('expression-suffix ','>, =, >=, <, <=, <>, λ',')'),
    @('simple-expr ','id, num, string, (, not, -',';, ), >, =, >=, <, <=, <>'),
    @('simple-expr-prime ','or, +, -, λ',';, ), >, =, >=, <, <=, <>'),
    @('term ','id, num, string, (, not, -','or, +, -'),
    @('term-prime ','and, *, /, λ','or, +, -'),
    @('factor-a ','id, num, string, (, not, -','and, *, /, or, +, -'),
    @('factor ','id, num, string, (','and, *, /, or, +, -'),
    @('relop ','>, =, >=, <, <=, <>','id, num, string, (, not, -'),
    @('addop ','or, +, -','id, num, string, (, not, -'),
    @('mulop ','and, *, /','id, num, string, (, not, -'),
    @('constant ','num, string','and, *, /, or, +, -')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# Remove leftover rows 34-36 from the old, longer table
$ws.Range("A34:C36").ClearContents()

# Update the sheet view: scrolled so row 14 is at top, active/selected cell C33
$ws.Range("C33").Select()
$excel.ActiveWindow.ScrollRow = 14
